$d = $word.ActiveDocument

# 1. Merge the "TEAM ID:" value runs: the separate " " run and the
#    "PNT2022TMID43416" run collapse into a single run reading
#    " PNT2022TMID43416" (same run formatting, since both runs already
#    shared identical rPr).
$d.Content.Find.Execute(" PNT2022TMID43416", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " PNT2022TMID43416", 2) | Out-Null

# 2. Remove the now-stale "_GoBack" bookmark pair.
$d.Bookmarks("_GoBack").Delete()
